$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("Prova1", "Prova2")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 3 is the template row; duplicate it into rows 4 and 5
    $srcRow = $ws.Range("A3:J3")

    foreach ($targetRow in 4, 5) {
        $dstRow = $ws.Range("A" + $targetRow + ":J" + $targetRow)
        $srcRow.Copy($dstRow)
    }
}
